# A new weekly price-record row for "Vega Monumental Concepción" / Cilantro
# needs to be inserted into the historical log, right above the existing
# row that is currently at row 295. All rows from 295 downward shift down
# by one (the very last existing row, 343, ends up at 344).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 295 - this shifts rows 295:343 down
# to 296:344 (and extends the sheet dimension to R344).
$ws.Range("A295").EntireRow.Insert()

# Populate the newly-inserted row 295 with the new record.
$ws.Range("A295").Value = 11
$ws.Range("B295").Value = 'Vega Monumental Concepción'
$ws.Range("C295").Value = 'Bíobío'
$ws.Range("D295").Value = 45154
$ws.Range("E295").Value = 8
$ws.Range("F295").Value = 100112040
$ws.Range("G295").Value = 'Cilantro'
$ws.Range("H295").Value = 'Sin especificar'
$ws.Range("I295").Value = 'Primera'
$ws.Range("J295").Value = 80
$ws.Range("K295").Value = 6500
$ws.Range("L295").Value = 6500
$ws.Range("M295").Value = 6500
$ws.Range("N295").Value = '$/caja 36 atados'
$ws.Range("O295").Value = 'Región Metropolitana'
$ws.Range("P295").Value = 181
$ws.Range("Q295").Value = 36
$ws.Range("R295").Value = 'Hortaliza'
